# Applies the "parte 1 de nuevos estado de cuenta" update:
#  - Removes the data row belonging to ELVIS ENRIQUE AVILA AVILA (CC 10904108)
#  - Updates VALOR MORA total, Cant. Trabajadores count, and the Salario Basico
#    for STEVEN MARTINEZ OLMOS's remaining row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Delete the entire worksheet row that holds the ELVIS ENRIQUE AVILA AVILA record.
# Deleting the row (rather than just clearing cells) shifts every row below it
# up by one, which is exactly what happened between the two workbook revisions
# (the signature block moved from rows 23-24 up to rows 22-23).
$ws.Rows("18:18").Delete()

# Update the three figures that changed on this revision.
$ws.Range("E11").Value = 226000     # VALOR MORA
$ws.Range("C13").Value = 2          # Cant. Trabajadores
$ws.Range("G17").Value = 3500000    # Salario Basico (STEVEN MARTINEZ OLMOS)
